$d = $word.ActiveDocument

# Plain substring replace-all (used for strings that are not at risk of
# colliding with other find/replace pairs in this batch).
function ReplaceAll($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

# Whole-word replace-all (used for short bare numbers, so we never touch a
# number embedded inside a longer one).
function ReplaceWord($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

# Company name (header block, two occurrences)
ReplaceAll "Твой Домик" "Парковка на высоте"

# INN (two occurrences)
ReplaceAll "ИНН: 423519337037" "ИНН: 624436165689"

# Company address (three occurrences)
ReplaceAll "Челябинская область, город Магнитогорск, пр. Карла Маркса 172" "Тверская область, город Солнечногорск, наб. Сталина, 42"

# Vehicle make / model / plate
ReplaceAll "BMW" "Citroen"
ReplaceWord "M3" "C3"
ReplaceAll "Е 666 КХ 777" "К 474 НК 725"

# Vehicle owner name and address
ReplaceAll "Скибко Анжелика Александровна" "Яковлева Арина Ивановна"
ReplaceAll "г. Магнитогорск, ул. Суворова, 126" "г. Москва, ул. Родниковая, 47, оф. 96"

# Dates / time (do the short bare numbers before the phone number, since the
# new phone number happens to contain "15" as a substring)
ReplaceAll "26.01.2025" "27.01.2025"
ReplaceWord "15" "21"
ReplaceWord "48" "43"
ReplaceAll "25.02.2025" "27.02.2025"

# Storekeeper name
ReplaceAll "Клементьева Анжелика Александровна" "Кривоносов Иван Алексеевич"

# Days stored / amount
ReplaceWord "30" "31"
ReplaceAll "5970" "6603"

# Owner phone number (last, see note above)
ReplaceAll "79514597925" "79395080159"
